$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Update D2 value from "Automation2" to "Automation3"
$ws.Range("D2").Value = "Automation3"

# Update the active selection to B3
$ws.Range("B3").Select() | Out-Null
